$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 185 (this shifts existing rows 185-205
#    down to 186-206, carrying their cell values/styles with them).
# ---------------------------------------------------------------------------
$ws.Rows.Item(185).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the newly inserted row 185 with the data for
#    "The Continent of International Law".
# ---------------------------------------------------------------------------
$ws.Cells.Item(185, 1).Value2 = "The Continent of International Law"
$ws.Cells.Item(185, 2).Value2 = "international relations"
$ws.Cells.Item(185, 3).Value2 = "http://www.isr.umich.edu/cps/coil/files.html"
$ws.Cells.Item(185, 4).Value2 = "International institutions, international law, international agreements, incomplete contracts"
$ws.Cells.Item(185, 5).Value2 = "world"
$ws.Cells.Item(185, 6).Value2 = 1925
$ws.Cells.Item(185, 7).Value2 = 1997
$ws.Cells.Item(185, 8).Value2 = "online"
$ws.Cells.Item(185, 9).Value2 = "free, no registration"
$ws.Cells.Item(185, 10).Value2 = "http://www.isr.umich.edu/cps/coil/COIL-instrument-100707.pdf"
$ws.Cells.Item(185, 11).Value2 = "http://www.isr.umich.edu/cps/coil/COIL-data-130123.csv"

# ---------------------------------------------------------------------------
# 3. Rebuild the worksheet hyperlinks.
#
#    Inserting the row above shifts cell contents down, but this COM
#    runtime does not shift the associated Hyperlinks collection
#    automatically, and individual Hyperlink.Delete() calls are inert here.
#    The only effective way to clear stale hyperlinks is a full
#    Cells.Hyperlinks.Delete(), so every hyperlink (existing, row-shifted,
#    and the 3 brand-new ones for row 185) is re-added from scratch below,
#    in the same relative order as the original workbook.
# ---------------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

$hyperlinkData = @(
    @('C7', 'http://www.electionstudies.org/', $null),
    @('C78', 'http://www.europeansocialsurvey.org/', $null),
    @('C8', 'http://www.vanderbilt.edu/lapop/about-americasbarometer.php', $null),
    @('C13', 'http://www.asianbarometer.org/', $null),
    @('C124', 'https://www.lissdata.nl/', $null),
    @('C149', 'http://www.pewglobal.org/datasets/', $null),
    @('C205', 'http://www.worldvaluessurvey.org/wvs.jsp', $null),
    @('C48', 'http://www.correlatesofwar.org/', $null),
    @('C41', 'http://www.cses.org/', $null),
    @('C69', 'http://ec.europa.eu/commfrontoffice/publicopinion/index.cfm', $null),
    @('C76', 'https://www.eurofound.europa.eu/surveys/european-quality-of-life-surveys', $null),
    @('C4', 'http://www.afrobarometer.org/', $null),
    @('C28', 'http://caucasusbarometer.org/en/datasets/', $null),
    @('C80', 'http://www.europeanvaluesstudy.eu/', $null),
    @('C113', 'http://www.issp.org/menu-top/home/', $null),
    @('C40', 'http://www.cpds-data.org/', $null),
    @('C144', 'http://www.edac.eu/policies_desc.cfm?v_id=112', $null),
    @('C57', 'https://sites.duke.edu/democracylinkage/data/', $null),
    @('C94', 'https://www.eui.eu/Research/Library/ResearchGuides/Economics/Statistics/DataPortal/GSOEP', $null),
    @('C24', 'https://www.iser.essex.ac.uk/bhps', $null),
    @('C192', 'https://www.understandingsociety.ac.uk/documentation/mainstage', $null),
    @('C10', 'http://arabbarometer.org/', $null),
    @('C52', 'http://www.valgprojektet.dk/default.asp?l=eng', $null),
    @('C73', 'http://europeanelectionstudies.net/ees-study-components/voter-study/', $null),
    @('C70', 'http://europeanelectionstudies.net/ees-study-components/euromanifesto-study/', $null),
    @('C74', 'http://europeanelectionstudies.net/ees-study-components/elite-study/', $null),
    @('C75', 'http://europeanelectionstudies.net/ees-study-components/media-study/', $null),
    @('C91', 'http://gss.norc.org/', $null),
    @('C121', 'http://www.latinobarometro.org/lat.jsp', $null),
    @('C179', 'http://forscenter.ch/en/our-surveys/swiss-household-panel/', $null),
    @('C127', 'https://manifesto-project.wzb.eu/datasets', $null),
    @('C107', 'http://www.impic-project.eu/data/', $null),
    @('C134', 'http://www.queensu.ca/mcp/about/definitionsdata', $null),
    @('C92', 'https://www.gesis.org/en/institute/competence-centers/rdc-allbus/', $null),
    @('C25', 'http://www.natcen.ac.uk/our-research/research/british-social-attitudes/', $null),
    @('C177', 'http://valforskning.pol.gu.se/english', $null),
    @('C23', 'http://www.britishelectionstudy.com/', $null),
    @('C202', 'http://www.ipu.org/wmn-e/world-arc.htm', $null),
    @('C111', 'https://havardhegre.net/iaep/', $null),
    @('C77', 'http://www.erdda.se/index.php/projects/erd/data-archive', $null),
    @('C156', 'http://www.erdda.se/index.php/projects/cpd/data-archive', $null),
    @('C26', 'http://www.erdda.se/index.php/projects/cpd/data-archive', $null),
    @('C71', 'http://sdw.ecb.europa.eu/', $null),
    @('C34', 'http://www.erdda.se/index.php/projects/cpd/data-archive', $null),
    @('C145', 'https://fsw.vu.nl/en/departments/political-science-and-public-administration/staff/woldendorp/party-government-data-set/index.aspx', $null),
    @('C199', 'https://sites.lsa.umich.edu/tsebelis/data/veto-players-data/', $null),
    @('C38', 'http://comparativeconstitutionsproject.org/', 'http://comparativeconstitutionsproject.org/'),
    @('C142', 'http://www.parlgov.org/', $null),
    @('C195', 'http://ucdp.uu.se/downloads/', $null),
    @('C160', 'http://qog.pol.gu.se/data/datadownloads/qogstandarddata', $null),
    @('C154', 'http://www.systemicpeace.org/polity/polity4.htm', $null),
    @('C32', 'https://www.chesdata.eu/our-surveys/', $null),
    @('C146', 'http://www.tcd.ie/Political_Science/ppmd/', $null),
    @('C67', 'https://mepsurvey.eu/data-objects/data/', $null),
    @('C100', 'http://www.start.umd.edu/gtd/', $null),
    @('C103', 'https://www.binghamton.edu/political-science/research.html/', $null),
    @('C143', 'https://www.binghamton.edu/political-science/research.html/', $null),
    @('C129', 'https://www.binghamton.edu/political-science/research.html/', $null),
    @('C166', 'https://www.binghamton.edu/political-science/research.html/', $null),
    @('C42', 'http://cwed2.org/download.php', $null),
    @('C181', 'http://web.missouri.edu/~williamslaro/mipdata.html', $null),
    @('C11', 'http://www.rochester.edu/college/faculty/hgoemans/data.htm', $null),
    @('C44', 'http://www.electiondataarchive.org/', $null),
    @('C184', 'http://www.humanrightsdata.com/p/data-documentation.html', $null),
    @('C22', 'https://sites.google.com/site/mkmtwo/data', $null),
    @('C197', 'https://www.prio.org/Data/Governance/Vanhanens-index-of-democracy/', $null),
    @('C36', 'http://www.comparativeagendas.net/datasets_codebooks', $null),
    @('C35', 'http://www.comparativeagendas.net/datasets_codebooks', $null),
    @('C37', 'http://www.comparativeagendas.net/datasets_codebooks', $null),
    @('C72', 'http://www.nsd.uib.no/european_election_database', $null),
    @('C200', 'https://www.idea.int/data-tools/data/voter-turnout', $null),
    @('C39', 'http://www.marquette.edu/polisci/faculty_swank.shtml', $null),
    @('C150', 'http://comparativepolitics.uni-greifswald.de/data.html', $null),
    @('C43', 'http://www.lisdatacenter.org/resources/other-databases/', $null),
    @('C63', 'https://www.gesis.org/angebot/daten-analysieren/weitere-sekundaerdaten/weitere-internationale-daten/europaeische-wahlstudien/election-studies-eastern-europe/', $null),
    @('C170', 'http://www.spin.su.se/datasets/scip', $null),
    @('C58', 'http://mattgolder.com/elections', $null),
    @('C53', 'http://econ.worldbank.org/WBSITE/EXTERNAL/EXTDEC/EXTRESEARCH/0,,contentMDK:20649465~pagePK:64214825~piPK:64214943~theSitePK:469382,00.html', $null),
    @('C88', 'http://www.frdb.org/page/data/scheda/frdb-iza-social-reforms-database/doc_pk/9027', $null),
    @('C6', 'http://aiddata.org/data/aiddata-core-research-release-level-1-v3-0', $null),
    @('C87', 'http://fundforpeace.org/fsi/data/', $null),
    @('C19', 'https://sites.google.com/site/authoritarianregimedataset/data', $null),
    @('C141', 'http://dx.doi.org/10.7910/DVN/ZTPW0Y', $null),
    @('C152', 'http://www.politicalterrorscale.org/Data/Download.html', $null),
    @('C187', 'https://mgmt.wharton.upenn.edu/faculty/heniszpolcon/polcondataset/', $null),
    @('C49', 'https://www.transparency.org/research/cpi/', $null),
    @('C198', 'https://www.v-dem.net/en/data/', $null),
    @('C62', 'https://www.fraserinstitute.org/economic-freedom/dataset', $null),
    @('C204', 'https://rsf.org/en/ranking_table', $null),
    @('C125', 'http://www.ggdc.net/maddison/maddison-project/data.htm', $null),
    @('C85', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C126', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C173', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C104', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C131', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C174', 'http://www.systemicpeace.org/inscrdata.html', $null),
    @('C81', 'http://ec.europa.eu/eurostat/data/bulkdownload', $null),
    @('C148', 'https://sites.google.com/site/electoralintegrityproject4/projects/expert-survey-2', $null),
    @('C64', 'http://hdl.handle.net/1902.1/17901', $null),
    @('C59', 'https://www.hertie-school.org/en/governancereport/govreport-indicators/', $null),
    @('C108', 'https://www.hertie-school.org/en/governancereport/govreport-indicators/', $null),
    @('C79', 'https://www.hertie-school.org/en/governancereport/govreport-indicators/', $null),
    @('C3', 'https://www.hertie-school.org/en/governancereport/govreport-indicators/', $null),
    @('C101', 'https://www.hertie-school.org/en/governancereport/govreport-indicators/', $null),
    @('C105', 'http://mo.ibrahim.foundation/iiag/downloads/', $null),
    @('C136', 'http://www.nelda.co/', $null),
    @('C66', 'http://epi.yale.edu/downloads', $null),
    @('C109', 'http://www.uva-aias.net/en/ictwss', $null),
    @('C172', 'http://www.brown.edu/Departments/Economics/Faculty/Louis_Putterman/antiquity%20index.htm', $null),
    @('C21', 'http://archive.ipu.org/gpr-e/downloads/index.htm', $null),
    @('C5', 'http://archive.ipu.org/gpr-e/downloads/index.htm', $null),
    @('C203', 'http://www.brown.edu/Departments/Economics/Faculty/Louis_Putterman/world%20migration%20matrix.htm', $null),
    @('C117', 'http://globalization.kof.ethz.ch/', $null),
    @('C110', 'https://sites.google.com/site/aljazkuncic/research', $null),
    @('C190', 'https://data.worldbank.org/data-catalog/worldwide-governance-indicators', $null),
    @('C165', 'https://worldjusticeproject.org/our-work/wjp-rule-law-index/wjp-rule-law-index-2016/current-historical-data', $null),
    @('C83', 'http://privatewww.essex.ac.uk/~ksg/exptradegdp.html', $null),
    @('C128', 'http://econ.worldbank.org/WBSITE/EXTERNAL/EXTDEC/EXTRESEARCH/0,,contentMDK:20699070~pagePK:64214825~piPK:64214943~theSitePK:469382,00.html', $null),
    @('C96', 'https://www.globalintegrity.org/downloads/', $null),
    @('C206', 'http://wid.world/data/', $null),
    @('C86', 'http://www.anderson.ucla.edu/faculty_pages/romain.wacziarg/papersum.html', $null),
    @('C163', 'http://www.thearda.com/archive/files/Descriptions/RASCONS.asp', $null),
    @('C135', 'https://unstats.un.org/unsd/snaama/dnlList.asp', $null),
    @('C102', 'https://www.wider.unu.edu/database/additional-data-files-grd', $null),
    @('C133', 'http://www.mar.umd.edu/mar_data.asp', $null),
    @('C171', 'https://www.strausscenter.org/scad.html', $null),
    @('C114', 'http://www.ines.tau.ac.il/elections.html', $null),
    @('C18', 'http://www.autnes.at/en/data-download/', $null),
    @('C27', 'http://ces-eec.arts.ubc.ca/english-section/surveys/', $null),
    @('C161', 'https://www.arjanschakel.nl/regauth_dat.html', $null),
    @('C95', 'http://www.globalelectionsdatabase.com/index.php/datasets', $null),
    @('C180', 'http://folk.uib.no/sspje/tweed.htm', $null),
    @('C158', 'https://www.wzb.eu/en/research/completed-research-programs/civil-society-and-political-mobilization/projects/prodat-dokumentation-un', $null),
    @('C157', 'https://www.unige.ch/sciences-societe/incite/welcome-to-the-incite-website/data/new-s/', $null),
    @('C159', 'https://github.com/erikgahner/polls', $null),
    @('C189', 'http://doi.org/10.7910/DVN/28856', $null),
    @('C178', 'https://github.com/MansMeg/SwedishPolls', $null),
    @('C196', 'https://github.com/zonination/election-history', $null),
    @('C51', 'https://github.com/Straubinger/folketingsvalg', $null),
    @('C116', 'http://faculty.tuck.dartmouth.edu/rafael-laporta/research-publications/', $null),
    @('J3', 'https://www.hertie-school.org/fileadmin/2_Research/1_About_our_research/4_The_Governance_Report/Indicators/2014/Downloads/GovReport2014_Indicators_IndicatorsCodebook.pdf', $null),
    @('M4', 'http://afrobarometer.org/sites/default/files/data/round-6/merged_r6_data_2016_36countries2.sav', $null),
    @('J4', 'http://afrobarometer.org/sites/default/files/data/round-6/merged_round_6_codebook_28082017.pdf', $null),
    @('N5', 'http://archive.ipu.org/gpr-e/downloads/data-age-gender-profession-e.xls', $null),
    @('J5', 'http://archive.ipu.org/gpr-e/downloads/dataset-notes-e.pdf', $null),
    @('P6', 'https://github.com/AidData-WM/public_datasets/releases/download/v3.0/AidDataCore_ResearchRelease_Level1_v3.0.zip', $null),
    @('J6', 'http://docs.aiddata.org/ad4/files/inline/readme.pdf', $null),
    @('C140', 'https://doi.org/10.18712/NSD-NSD2405-V1', $null),
    @('C46', 'https://cces.gov.harvard.edu/', $null),
    @('P8', 'http://datasets.americasbarometer.org/database/files/746278534AmericasBarometer%20Grand%20Merge%202004-2014%20v3.0_FREE_dta.zip', $null),
    @('J8', 'http://datasets.americasbarometer.org/database/files/12364388022004-2014%20Grand%20Merge%20Codebook_V3.0_Free_W.pdf', $null),
    @('M10', 'http://www.arabbarometer.org/sites/default/files/Arab_Barometer_Fourth_Wave_English_Data_Set_v1.sav', $null),
    @('L10', 'http://www.arabbarometer.org/sites/default/files/Arab_Barometer_Fourth_Wave_English_Data_Set_v1.dta', $null),
    @('J10', 'http://www.arabbarometer.org/sites/default/files/code_book/AB4%20English%20Codebook%20Final%20English.pdf', $null),
    @('C194', 'http://www.uni-heidelberg.de/fakultaeten/wiso/awi/professuren/intwipol/datasets_en.html', $null),
    @('N194', 'http://www.axel-dreher.de/UNSCdata.xls', $null),
    @('C29', 'https://www.kof.ethz.ch/services/daten/data-on-central-bank-governors.html', $null),
    @('N29', 'https://www.ethz.ch/content/dam/ethz/special-interest/dual/kof-dam/documents/central_bank_governors/cbg_turnover.xlsx', $null),
    @('C82', 'http://dx.doi.org/10.7910/DVN/XPCVEI', $null),
    @('C188', 'https://snd.gu.se/en/catalogue/study/snd0905', $null),
    @('J188', 'https://snd.gu.se/catalogue/file/5665', $null),
    @('C186', 'https://drryanmaness.wixsite.com/cyberconflcit/cyber-conflict-dataset', $null),
    @('N186', 'https://docs.wixstatic.com/ugd/4b99a4_294fde43b8094872999ca63f62972765.xlsx?dn=DCID%20Version%201.1.xlsx', $null),
    @('J186', 'https://docs.wixstatic.com/ugd/4b99a4_4c7971ea7791464a8ac551fff85fb1f1.pdf', $null),
    @('C122', 'http://dx.doi.org/10.7910/DVN/SYZZEY', $null),
    @('C120', 'http://dx.doi.org/10.7910/DVN/24872', $null),
    @('C54', 'http://www.robertthomson.info/research/resolving-controversy-in-the-eu', $null),
    @('J54', 'http://www.robertthomson.info/wp-content/uploads/2011/01/Issues_list_new_26March2012.pdf', $null),
    @('K54', 'http://www.robertthomson.info/wp-content/uploads/2011/01/deu15_27_26March2012.csv', $null),
    @('C130', 'https://zenodo.org/record/61234', $null),
    @('N130', 'https://zenodo.org/record/61234/files/MAPP_dataset_-_Version_2.0.xlsx', $null),
    @('C97', 'http://faculty.uml.edu/Jenifer_whittenwoodring/MediaFreedomData_000.aspx', $null),
    @('K97', 'http://faculty.uml.edu/Jenifer_whittenwoodring/GMFD_V2.csv', $null),
    @('C201', 'http://mediaproject.wesleyan.edu/dataaccess/', $null),
    @('C153', 'http://politicaladarchive.org/data/', $null),
    @('K153', 'http://politicaladarchive.org/api/v1/ad_instances?output=csv', $null),
    @('L11', 'http://www.rochester.edu/college/faculty/hgoemans/Archigos_4.1_stata14.dta', $null),
    @('J11', 'http://www.rochester.edu/college/faculty/hgoemans/Archigos_4.1.pdf', $null),
    @('J13', 'http://www.asianbarometer.org/pdf/core_questionnaire_wave4.pdf', $null),
    @('C167', 'http://faculty.missouri.edu/williamslaro/govtdata.html', $null),
    @('J167', 'http://faculty.missouri.edu/williamslaro/SW%202016%20Codebook--Governments.pdf', $null),
    @('L167', 'http://faculty.missouri.edu/williamslaro/Seki-Williams%20Governments--Version%202.0.dta', $null),
    @('K167', 'http://faculty.missouri.edu/williamslaro/Seki-Williams%20Governments--Version%202.0.csv', $null),
    @('C168', 'http://faculty.missouri.edu/williamslaro/govtdata.html', $null),
    @('J168', 'http://faculty.missouri.edu/williamslaro/SW%202016%20Codebook--Ministers.pdf', $null),
    @('K168', 'http://faculty.missouri.edu/williamslaro/Seki-Williams%20Ministers--Version%202.0.csv', $null),
    @('L168', 'http://faculty.missouri.edu/williamslaro/Seki-Williams%20Ministers--Version%202.0.dta', $null),
    @('C84', 'https://sites.google.com/site/econolaols/extended-state-history-index', $null),
    @('J84', 'https://drive.google.com/file/d/1t5p1USIivXK-38urc2d5Fx7X5rHTzxzQ/view?usp=sharing', $null),
    @('N84', 'https://drive.google.com/file/d/1OExJ70t0YBzlV1vhUE11bbQTW8tJH7cW/view?usp=sharing', $null),
    @('C176', 'http://www.share-project.org/data-documentation/share-data-releases.html', $null),
    @('C123', 'http://www.ebrd.com/what-we-do/economic-research-and-data/data/lits.html', $null),
    @('J123', 'http://www.ebrd.com/cs/Satellite?c=Content&cid=1395256887325&d=&pagename=EBRD%2FContent%2FDownloadDocument', $null),
    @('L123', 'http://www.ebrd.com/cs/Satellite?c=Content&cid=1395256887465&d=&pagename=EBRD%2FContent%2FDownloadDocument', $null),
    @('C93', 'https://www.gesis.org/en/elections-home/gles/data/', $null),
    @('C68', 'http://www.columbia.edu/~aw2951/Datasets.html', $null),
    @('L68', 'http://www.columbia.edu/~aw2951/EPR3CountryNewReduced.dta', $null),
    @('K68', 'http://www.columbia.edu/~aw2951/eprnew301.xlsx', $null),
    @('J68', 'http://www.columbia.edu/~aw2951/AppendixEthnicPolitics.pdf', $null),
    @('C90', 'http://www.columbia.edu/~aw2951/Datasets.html', $null),
    @('J90', 'http://www.columbia.edu/~aw2951/EmpireNSdataset.pdf', $null),
    @('K90', 'http://www.columbia.edu/~aw2951/WimmerMin1.0.xls', $null),
    @('L90', 'http://www.columbia.edu/~aw2951/WimmerMin1.0.dta', $null),
    @('C155', 'https://www.markpack.org.uk/opinion-polls/', $null),
    @('N155', 'https://www.markpack.org.uk/files/2017/10/PollBase.xls', $null),
    @('C50', 'https://www.ctdatacollaborative.org/download-global-dataset', $null),
    @('K50', 'https://www.ctdatacollaborative.org/sites/default/files/The%20Global%20Dataset%2027%20Nov%202017_0.csv', $null),
    @('J50', 'https://www.ctdatacollaborative.org/sites/default/files/CTDC%20codebook%20v6_0.pdf', $null),
    @('C12', 'https://www.acleddata.com/data/acled-version-7-1997-2016/', $null),
    @('J12', 'https://www.acleddata.com/wp-content/uploads/2017/01/ACLED_Codebook_2017.pdf', $null),
    @('P12', 'https://www.acleddata.com/wp-content/uploads/2017/01/ACLED-Version-7-All-Africa-1997-2016_csv_dyadic-file.zip', $null),
    @('N12', 'https://www.acleddata.com/wp-content/uploads/2017/01/ACLED-Version-7-All-Africa-1997-2016_dyadic-file.xlsx', $null),
    @('C56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Data.htm', $null),
    @('J56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Democracy%20CrossNational%20Data/Democracy%20Crossnational%20Codebook%20March%202009.pdf', $null),
    @('K56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Democracy%20CrossNational%20Data/Democracy%20Crossnational%20Data%20Spring%202009.csv', $null),
    @('L56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Democracy%20CrossNational%20Data/Democracy%20Crossnational%20Data%20Spring%202009%20StataSE.dta', $null),
    @('M56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Democracy%20CrossNational%20Data/Democracy%20Crossnational%20Data%20Spring%202009.sav', $null),
    @('N56', 'https://sites.hks.harvard.edu/fs/pnorris/Data/Democracy%20CrossNational%20Data/Democracy%20Crossnational%20Data%20Spring%202009%20Excel.xlsx', $null),
    @('C30', 'http://dx.doi.org/10.7910/DVN/ALVXLM', $null),
    @('C193', 'http://hdl.handle.net/1902.1/12379', $null),
    @('C47', 'http://www.ippsr.msu.edu/public-policy/correlates-state-policy', $null),
    @('J47', 'http://ippsr.msu.edu/sites/default/files/CorrelatesCodebook.pdf', $null),
    @('N47', 'http://ippsr.msu.edu/sites/default/files/correlatesofstatepolicyprojectv1_13.xlsx', $null),
    @('L47', 'http://ippsr.msu.edu/sites/default/files/correlatesofstatepolicyprojectv1_13.dta', $null),
    @('K47', 'http://ippsr.msu.edu/sites/default/files/correlatesofstatepolicyprojectv1_13.csv', $null),
    @('C191', 'https://www.prio.org/JPR/Datasets/', $null),
    @('P191', 'http://file.prio.no/Journals/JPR/2018/55/1/Sara%20Lindberg%20Bromley.zip', $null),
    @('C132', 'http://dx.doi.org/10.7910/DVN/FB0R8A', $null),
    @('C151', 'http://journals.sagepub.com/doi/abs/10.1177/0022343316628932', $null),
    @('P151', 'http://file.prio.no/journals/JPR/2016/53/4/Nadine%20Ansorg,%20Felix%20Haass%20&%20Julia%20Strasheim.zip', $null),
    @('C89', 'http://www.csae.ox.ac.uk/general/free-and-fair-elections-database', $null),
    @('J89', 'http://www.csae.ox.ac.uk/materials/data/151/csae-data-ffelectionscodebookmarch2014.pdf', $null),
    @('K89', 'http://www.csae.ox.ac.uk/materials/data/151/csae-data-ffelections-v11-dates.csv', $null),
    @('N89', 'http://www.csae.ox.ac.uk/materials/data/151/csae-data-ffelections-v11-dates.xlsx', $null),
    @('C164', 'http://www.efetokdemir.com/data.html', $null),
    @('J164', 'http://www.efetokdemir.com/uploads/3/7/3/2/37326045/rtg_dataset_-_codebook.pdf', $null),
    @('L164', 'http://www.efetokdemir.com/uploads/3/7/3/2/37326045/replicationdatajpr-oldstata.dta', $null),
    @('C45', 'http://dx.doi.org/10.7910/DVN/F8ITEB', $null),
    @('C31', 'http://www.chisols.org/', $null),
    @('J31', 'http://www.chisols.org/uploads/1/1/2/6/11264284/chisolsusermanualv4.0.pdf', $null),
    @('P31', 'http://www.chisols.org/uploads/1/1/2/6/11264284/chisolsstyr4_0.zip', $null),
    @('C147', 'https://peaceaccords.nd.edu/research', $null),
    @('N147', 'http://peaceaccords.nd.edu/sites/default/files/PAM_ID%20V.1.5%20Updated%2029JULY2015.xlsx', $null),
    @('J147', 'http://peaceaccords.nd.edu/sites/default/files/PAM_ID%20CODEBOOK%20V.1.5%2029July2015.pdf', $null),
    @('C183', 'http://americanideologyproject.com/', $null),
    @('J183', 'http://americanideologyproject.com/estimates/estimates2015/codebook.pdf', $null),
    @('C139', 'https://www.du.edu/korbel/sie/research/chenow_navco_data.html', $null),
    @('C20', 'http://sites.psu.edu/dictators/', $null),
    @('J20', 'http://sites.psu.edu/dictators/wp-content/uploads/sites/12570/2016/05/GWF-Codebook.pdf', $null),
    @('P20', 'http://sites.psu.edu/dictators/wp-content/uploads/sites/12570/2016/05/GWF-Autocratic-Regimes-1.2.zip', $null),
    @('C169', 'http://www.sexualviolencedata.org/dataset/', $null),
    @('J169', 'http://www.sexualviolencedata.org/wp-content/uploads/2013/01/SVAC-coding-manual-10-25-13.pdf', $null),
    @('P169', 'http://www.sexualviolencedata.org/wp-content/uploads/2013/01/SVAC_dataset-update-2016-June-21.xlsx.zip', $null),
    @('C16', 'http://www.australianelectionstudy.org/voter_studies.html', $null),
    @('C14', 'http://www.australianelectionstudy.org/candidate_studies.html', $null),
    @('C15', 'http://www.australianelectionstudy.org/other_studies_acrs.html', $null),
    @('C17', 'http://www.australianelectionstudy.org/anpas.html', $null),
    @('M15', 'http://ada.edu.au/ADAData/data/acrs_1999_01018-pub.sav', $null),
    @('J15', 'http://nesstar.ada.edu.au/webview/velocity?study=http://150.203.254.120:80/obj/fStudy/au.edu.anu.ada.ddi.01018-pub&format=pdf&mode=transform&s&gs', $null),
    @('M16', 'http://ada.edu.au/ADAData/data/aes_2016_01365.sav', $null),
    @('J16', 'http://nesstar.ada.edu.au/webview/velocity?study=http://150.203.254.120:80/obj/fStudy/au.edu.anu.ada.ddi.01365&format=pdf&mode=transform&s&gs', $null),
    @('M14', 'http://ada.edu.au/ADAData/data/acs_2016_01366.sav', $null),
    @('J14', 'http://ada.edu.au/ADAData/AES/Australian%20Candidate%20Study%202016.pdf', $null),
    @('M17', 'http://ada.edu.au/ADAData/data/anpas_1979_00009.sav', $null),
    @('J17', 'http://nesstar.ada.edu.au/webview/velocity?study=http://150.203.254.120:80/obj/fStudy/au.edu.anu.ada.ddi.00009&format=pdf&mode=transform&s&gs', $null),
    @('J100', 'http://www.start.umd.edu/gtd/downloads/Codebook.pdf', $null),
    @('C175', 'http://cpostdata.uchicago.edu/search_new.php', $null),
    @('C55', 'http://www.democracybarometer.org/dataset_en.html', $null),
    @('J55', 'http://www.democracybarometer.org/Data/Codebook_all%20countries_1990-2014.pdf', $null),
    @('L55', 'http://www.democracybarometer.org/Data/DB_data_1990-2014_Standardized.dta', $null),
    @('N55', 'http://www.democracybarometer.org/Data/DB_data_1990-2014_Standardized.xlsx', $null),
    @('C61', 'http://dx.doi.org/10.7910/DVN/UXIBNO', $null),
    @('C2', 'http://dx.doi.org/10.7910/DVN/29106', $null),
    @('N2', 'http://ps.au.dk/fileadmin/Statskundskab/Dokumenter/Forskning/Forskningscentre/DEDERE/lied_v3.xls', $null),
    @('C33', 'http://ps.au.dk/forskning/forskningsprojekter/dedere/datasets/', $null),
    @('J33', 'http://ps.au.dk/fileadmin/Statskundskab/Dokumenter/Forskning/Forskningscentre/DEDERE/CLDcodebook.pdf', $null),
    @('N33', 'http://ps.au.dk/fileadmin/Statskundskab/Dokumenter/Forskning/Forskningscentre/DEDERE/CLD_cow.xls', $null),
    @('C106', 'http://faculty.ucmerced.edu/cconrad2/Academic/Data.html', $null),
    @('P106', 'http://faculty.ucmerced.edu/cconrad2/Academic/Data_files/CY.zip', $null),
    @('J106', 'http://faculty.ucmerced.edu/cconrad2/Academic/Data_files/ITT_CY_UsersGuide19July11.pdf', $null),
    @('C118', 'https://cms.uni-konstanz.de/fileadmin/archive/kosved/polver/gschneider/forschung/kosved/data/', $null),
    @('J118', 'https://cms.uni-konstanz.de/fileadmin/polver/gschneider/KOSVED/Coders_and_sources_kosved_webappendix_20130527.pdf', $null),
    @('N118', 'https://cms.uni-konstanz.de/fileadmin/polver/gschneider/KOSVED/KOSVED_Gesamtdatensatz_20130529_ohneGIS.xlsx', $null),
    @('C60', 'http://hdl.handle.net/1902.1/14717', $null),
    @('C137', 'https://cise.luiss.it/cise/dataset-of-new-parties-and-party-system-innovation-in-western-europe-since-1945/', $null),
    @('N137', 'https://cise.luiss.it/cise/wp-content/uploads/downloads/2017/08/Dataset-of-New-Parties-and-Party-System-Innovation-after-1945.xlsx', $null),
    @('J137', 'https://cise.luiss.it/cise/wp-content/uploads/downloads/2017/08/New-parties-and-party-system-innovation-codebook-and-dataset-information.pdf', $null),
    @('J180', 'http://folk.uib.no/sspje/TWEED%20Code%20Book.pdf', $null),
    @('P180', 'http://folk.uib.no/sspje/tweed.zip', $null),
    @('C65', 'https://cise.luiss.it/cise/dataset-of-electoral-volatility-and-its-internal-components-in-western-europe-1945-2015/', $null),
    @('J65', 'https://cise.luiss.it/cise/download/Codebook-and-dataset-information1.pdf', $null),
    @('N65', 'https://cise.luiss.it/cise/wp-content/uploads/downloads/2017/08/Dataset-of-Electoral-Volatility-and-its-internal-components-in-Western-Europe-1945-2015.xlsx', $null),
    @('C99', 'https://sites.google.com/a/thomaserichter.de/gsre/gsre-1-0', $null),
    @('C98', 'https://www.idea.int/gsod-indices/dataset-resources', $null),
    @('J98', 'https://www.idea.int/gsod/files/IDEA-GSOD-2017-CODEBOOK.pdf', $null),
    @('K98', 'https://www.idea.int/gsod-indices/sites/default/files/gsodi_pv_1.csv', $null),
    @('N98', 'https://www.idea.int/gsod-indices/sites/default/files/gsodi_pv_1.xlsx', $null),
    @('M98', 'https://www.idea.int/gsod-indices/sites/default/files/gsodi_pv_1.sav', $null),
    @('C119', 'https://sites.google.com/site/knemoto1978/kuniakinemotodata', $null),
    @('C138', 'https://sites.google.com/site/knemoto1978/kuniakinemotodata', $null),
    @('N119', 'https://sites.google.com/site/knemoto1978/NA%20Elections%201988-2016.xlsx', $null),
    @('N138', 'https://sites.google.com/site/knemoto1978/MMP%20Elections%201996-2014.xlsx', $null),
    @('C112', 'http://dx.doi.org/10.7910/DVN/X093TV', $null),
    @('C115', 'http://www.macrohistory.net/data/', $null),
    @('J115', 'http://www.macrohistory.net/JST/JSTdocumentationR2.pdf', $null),
    @('L115', 'http://www.macrohistory.net/JST/JSTdatasetR2.dta', $null),
    @('N115', 'http://www.macrohistory.net/JST/JSTdatasetR2.xlsx', $null),
    @('C9', 'http://folk.uio.no/bjornkho/MEP/', $null),
    @('C182', 'http://pages.ucsd.edu/~egartzke/datasets.htm', $null),
    @('J182', 'http://pages.ucsd.edu/~egartzke/data/affinity_codebook_03102006.pdf', $null),
    @('L182', 'http://pages.ucsd.edu/~egartzke/data/affinity_03102006.dta', $null),
    @('P182', 'http://pages.ucsd.edu/~egartzke/data/affinity_03102006.zip', $null),
    @('C162', 'http://hdl.handle.net/1902.1/16845', $null),
    @('C185', 'http://www.isr.umich.edu/cps/coil/files.html', $null),
    @('K185', 'http://www.isr.umich.edu/cps/coil/COIL-data-130123.csv', $null),
    @('J185', 'http://www.isr.umich.edu/cps/coil/COIL-instrument-100707.pdf', $null)
)

foreach ($item in $hyperlinkData) {
    $ref = $item[0]
    $target = $item[1]
    $display = $item[2]
    $rng = $ws.Range($ref)
    if ($display) {
        $ws.Hyperlinks.Add($rng, $target, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $display) | Out-Null
    } else {
        $ws.Hyperlinks.Add($rng, $target) | Out-Null
    }
}

# Re-adding a hyperlink re-applies the built-in "Hyperlink" character style to
# its cell, which can introduce a redundant duplicate of the existing
# "Hyperlink" cell style. Re-asserting the named style on every hyperlinked
# cell collapses it back onto the single pre-existing style actually used
# throughout the sheet.
foreach ($item in $hyperlinkData) {
    $ws.Range($item[0]).Style = "Hyperlink"
}

# ---------------------------------------------------------------------------
# 4. Refresh the sortState range so it covers the new last row (A2:Q206).
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:Q206"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 5. Restore the view/selection state to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("A185").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 174
$win.ScrollColumn = 1
